$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. NumberFormat is forced to
# Text ("@") before the write so numeric-looking strings (e.g. "215.90")
# are not silently coerced into actual numbers by Excel, then the style
# is reset back to "Normal" afterwards so no stray per-cell style index
# is left behind (matches the source workbook, where these cells carry
# no explicit style).
$updates = @{
    "D2" = "27.040.72"
    "E2" = "  +0.44%  "
    "D3" = "1.680.76"
    "E3" = "  +0.95%  "
    "E4" = "  +0.01%  "
    "D5" = "215.90"
    "E5" = "  +0.16%  "
    "E6" = "  -2.45%  "
    "E7" = "  -0.04%  "
    "E8" = "  +2.20%  "
    "D9" = "21.46"
    "E9" = "  +5.95%  "
    "D10" = "0.0625"
    "E10" = "  +0.75%  "
    "E11" = "  -0.88%  "
    "D12" = "1.917.51"
    "E12" = "  +0.91%  "
    "D13" = "1.704.11"
    "E13" = "  +2.35%  "
    "E14" = "  +0.93%  "
    "D15" = "0.533"
    "E15" = "  +1.72%  "
    "D16" = "66.51"
    "D17" = "27.033.55"
    "E17" = "  +0.45%  "
    "E18" = "  +2.08%  "
    "D19" = "236.13"
    "E19" = "  +0.80%  "
    "E20" = "  +0.85%  "
    "E21" = "  -0.01%  "
    "D22" = "4.48"
    "E22" = "  +2.94%  "
    "E23" = "  +1.89%  "
    "E24" = "  -3.83%  "
    "D25" = "146.60"
    "E25" = "  +0.28%  "
    "D27" = "16.49"
    "E27" = "  +3.79%  "
    "E29" = "  +0.00%  "
    "E30" = "  +0.38%  "
    "E31" = "  +0.19%  "
    "D33" = "1.537.40"
    "E33" = "  +5.53%  "
    "E34" = "  +2.22%  "
    "D35" = "1.73"
    "E35" = "  +5.89%  "
    "E36" = "  -0.86%  "
    "E37" = "  +1.52%  "
    "D38" = "0.919"
    "E38" = "  +1.68%  "
    "E39" = "  +3.22%  "
    "E40" = "  +6.68%  "
    "E41" = "  -0.01%  "
    "D42" = "67.96"
    "D43" = "5.59"
    "E43" = "  -2.00%  "
    "E44" = "  -0.41%  "
    "D45" = "1.822.60"
    "E45" = "  +0.49%  "
    "D46" = "0.779"
    "D47" = "90.52"
    "B48" = "RenderToken"
    "C48" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D48" = "1.54"
    "E48" = "  +0.29%  "
    "B49" = "Algorand"
    "C49" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "D49" = "0.104"
    "E49" = "  +2.62%  "
    "B50" = "EnergySwap"
    "C50" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D50" = "8.02"
    "E50" = "  +6.14%  "
    "B51" = "Cronos"
    "C51" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D51" = "0.0507"
    "E51" = "  +0.14%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
